$d = $word.ActiveDocument

# --- Edit 1: merge "Acessar_Dados" + "_Pessoais" runs into one run and drop the
#     now-redundant gramStart/gramEnd proofErr markers around them ---
$p1 = $d.Paragraphs(9)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="646FCFDE" w14:textId="680E3E0F" w:rsidR="003852BC" w:rsidRDefault="003852BC" w:rsidP="003852BC"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Read</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Acessar_Dados_Pessoais</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p1.Range.InsertXML($xml1)

# --- Edit 2: merge "Verificar_Status_de" + "_Chamad" runs into one run (keeping
#     the trailing "o" run separate) and drop the gramStart/gramEnd proofErr ---
$p2 = $d.Paragraphs(10)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7648523B" w14:textId="26B73297" w:rsidR="00A72556" w:rsidRDefault="003852BC" w:rsidP="00A72556"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>&#8211; Upload/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Verificar_Status_de_Chamad</w:t></w:r><w:r w:rsidR="00A72556"><w:t>o</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($xml2)

# --- Edit 3: remove the last 3 blank paragraphs before "Ultimo teste" and the
#     "Ultimo teste" paragraph itself (4 paragraphs total) ---
$startPos = $d.Paragraphs(16).Range.Start
$endPos = $d.Paragraphs(19).Range.End
$delRange = $d.Range($startPos, $endPos)
$delRange.Delete()
